$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

$ws.Range("B2").Value = "NO_WAIP_001"
$ws.Range("B3").Value = "AIP (acidification index periphyton)"
$ws.Range("B4").Value = "Norway"
$ws.Range("B5").Value = "Europe"
$ws.Range("B6").Value = "A2 - Chemical State characteristics"
$ws.Range("B7").Value = "Freshwater (F)"
$ws.Range("B8").Value = "F1 Rivers and streams biome"
$ws.Range("B9").Value = "F1.3 Freeze-thaw streams"
$ws.Range("B10").Value = 2025
$ws.Range("B11").Value = 2025
$ws.Range("B14").Value = "First draft"
$ws.Range("B15").Value = "Sandvik, H."
$ws.Range("B16").Value = "https://github.com/NINAnor/ecRxiv/tree/main/indicators/NO_WAIP_001"
$ws.Range("B17").Value = "Yes"
$ws.Range("B18").Value = "Scale - Transform - Aggregate (length-weighted arithmetic mean) - Truncate"

$ws.Activate()
$ws.Range("B2").Select()
